$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing style of the Price column, then temporarily force
# a text number format so numeric-looking strings (e.g. "1.00", "581.17")
# are stored as text, matching the original inline-string cell type.
$priceRange = $ws.Range("D2:D51")
$origStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '64.134.93'
$ws.Range('E2').Value = '  +5.30%  '
$ws.Range('D3').Value = '2.732.02'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '581.17'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '157.88'
$ws.Range('E6').Value = '  +9.56%  '
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  +4.89%  '
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '2.758.85'
$ws.Range('E9').Value = '  +3.97%  '
$ws.Range('E10').Value = '  +2.98%  '
$ws.Range('E11').Value = '  +3.86%  '
$ws.Range('E12').Value = '  +4.65%  '
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '3.220.24'
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').Value = '27.31'
$ws.Range('E15').Value = '  +4.18%  '
$ws.Range('D16').Value = '64.082.08'
$ws.Range('E16').Value = '  +5.27%  '
$ws.Range('E17').Value = '  +7.41%  '
$ws.Range('D18').Value = '2.760.30'
$ws.Range('E18').Value = '  +4.23%  '
$ws.Range('D19').Value = '12.08'
$ws.Range('E19').Value = '  +4.03%  '
$ws.Range('E20').Value = '  +4.66%  '
$ws.Range('D21').Value = '364.08'
$ws.Range('E21').Value = '  +3.49%  '
$ws.Range('D22').Value = '6.97'
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').Value = '0.548'
$ws.Range('E23').Value = '  +4.21%  '
$ws.Range('D24').Value = '0.997'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').Value = '67.14'
$ws.Range('E25').Value = '  +4.97%  '
$ws.Range('D26').Value = '0.173'
$ws.Range('E26').Value = '  +6.31%  '
$ws.Range('D27').Value = '8.65'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').Value = '0.0₃0925'
$ws.Range('E29').Value = '  +14.29%  '
$ws.Range('D30').Value = '2.02'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').Value = '7.24'
$ws.Range('D32').Value = '1.27'
$ws.Range('E32').Value = '  +18.22%  '
$ws.Range('D33').Value = '173.77'
$ws.Range('E33').Value = '  +3.84%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').Value = '  +3.54%  '
$ws.Range('D36').Value = '4.93'
$ws.Range('E36').Value = '  +7.40%  '
$ws.Range('E37').Value = '  +11.32%  '
$ws.Range('E38').Value = '  +8.88%  '
$ws.Range('E39').Value = '  +12.27%  '
$ws.Range('E40').Value = '  +4.83%  '
$ws.Range('D41').Value = '338.08'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('D42').Value = '6.12'
$ws.Range('E42').Value = '  +15.62%  '
$ws.Range('D43').Value = '39.45'
$ws.Range('E43').Value = '  +3.10%  '
$ws.Range('D44').Value = '21.93'
$ws.Range('E44').Value = '  +8.17%  '
$ws.Range('D45').Value = '22.26'
$ws.Range('E45').Value = '  +5.42%  '
$ws.Range('D46').Value = '0.0605'
$ws.Range('E46').Value = '  +5.58%  '
$ws.Range('E47').Value = '  +3.62%  '
$ws.Range('E48').Value = '  +4.60%  '
$ws.Range('D49').Value = '137.59'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('E51').Value = '  -0.01%  '

# Restore the original style/number format for the column.
$priceRange.Style = $origStyle
